# [PHOENIX-5849] Modified the data according to the production dump
#
# This script updates the financialTestData workbook so that the
# ledgerCode ("G Basheer Ahmed") is replaced with "KMC001" everywhere it
# is used, and the accountCode2 value used by the
# budgetCheckWithOutSubledger / voucherWithOutSubledger rows is changed
# from 3501003 to 1301001.

$wb = $excel.ActiveWorkbook

$wsJournal = $wb.Worksheets.Item("journalVoucherDetails")
$wsDirect  = $wb.Worksheets.Item("directBankPaymentDetails")

# ---------------------------------------------------------------------
# Replace the ledgerCode value "G Basheer Ahmed" with "KMC001" in every
# cell that references it (journalVoucherDetails + directBankPaymentDetails)
# ---------------------------------------------------------------------
$ledgerCodeCells = @("P2","Q2","P3","Q3","P4","Q4","P5","Q5","P6","P9","Q9")
foreach ($cellRef in $ledgerCodeCells) {
    $wsJournal.Range($cellRef).Value = "KMC001"
}

$wsDirect.Range("L2").Value = "KMC001"

# ---------------------------------------------------------------------
# Update accountCode2 (column G) for the budgetCheckWithOutSubledger (row 7)
# and voucherWithOutSubledger (row 8) rows from 3501003 to 1301001
# ---------------------------------------------------------------------
$wsJournal.Range("G7").Value = "1301001"
$wsJournal.Range("G8").Value = "1301001"
